$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B14").Value = 10
$ws.Range("E14").Value = "Screw Driver"
$ws.Range("E15").Value = "Monitor"

$ws.Range("E15:H15").Select()
